$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 168: 23.03.2024
$ws.Range("A168").Value = "23.03.2024"
$ws.Range("B168").Value = "22.03.2024"
$ws.Range("C168").Value = 32070
$ws.Range("D168").Value = 13000
$ws.Range("E168").Value = 8400
$ws.Range("F168").Value = 74298
$ws.Range("G168").Value = 8663
$ws.Range("H168").Value = 6327
$ws.Range("I168").Value = 8000
$ws.Range("J168").Value = 442
$ws.Range("K168").Value = 116
$ws.Range("L168").Value = 4650
$ws.Range("M168").Value = "https://web.archive.org/web/20240323053658/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker"

# New row 169: 24.03.2024
$ws.Range("A169").Value = "24.03.2024"
$ws.Range("B169").Value = "24.03.2024"
$ws.Range("C169").Value = 32226
$ws.Range("D169").Value = 13000
$ws.Range("E169").Value = 8400
$ws.Range("F169").Value = 74518
$ws.Range("G169").Value = 8663
$ws.Range("H169").Value = 6327
$ws.Range("I169").Value = 8000
$ws.Range("J169").Value = 442
$ws.Range("K169").Value = 116
$ws.Range("L169").Value = 4700
$ws.Range("M169").Value = "https://web.archive.org/web/20240324191339/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker"

# Row 169's F cell gets a wrap-text style in the target workbook
$ws.Range("F169").WrapText = $true

$ws.Range("M169").Select()
